$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBPPRTY")

# Relabel the header in A1: this row previously held "Max Potential Capacity (MW)"
# and now reflects the new retrofit-eligibility metric.
$ws.Range("A1").Value = "% of capacity available for retrofit"

# Row 19 is "hard coal w CCS". We no longer ban CCUS retrofits outright for all
# future years -- only years 2021-2027 (cols B:H) keep the ban flag of 1.
# Starting 2028 (col I) through 2050 (col AE), the ban flag is now 0, formatted
# as a plain integer (matching the other boolean rows' number format).
$retrofitRange = $ws.Range("I19:AE19")
$retrofitRange.Value = 0
$retrofitRange.NumberFormat = "0"

# Move the saved selection to A2 on the BBPPRTY sheet.
$ws.Range("A2").Select() | Out-Null
